# rill-analysis-report-core: Report params config
#
# - Point the "_settings" sheet's REST endpoint cell at the new
#   rill-analysis-web service (only the displayed text changes - the
#   underlying hyperlink relationship target is left as-is).
# - Widen column B on "_settings" so the longer URL fits.
# - Make "_settings" the active/selected sheet (it was "_input").

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("_settings")
$wsInput = $wb.Worksheets.Item("_input")

$wsSettings.Range("B1").Value = "http://10.81.21.140:8280/rill-analysis-web/rest/"

$wsSettings.Columns.Item(2).ColumnWidth = 53.15

$wsSettings.Activate()
